$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aVals = @(5440,5400,5370,5340,5320,5300,5280,5270,5260,5260,5270,5280,5290,5310,5340,5370,5430,5510,5610,5740,5880,6030,6200,6370,6530,6690,6830,6960,7070,7150,7220,7260,7280,7280,7270,7240,7210,7170,7130,7080,7030,6980,6930,6870,6810,6750,6690,6630,6570,6520,6480,6430,6400,6380,6360,6350,6350,6360,6370,6410,6480,6540,6620,6700,6790,6880,6970,7050,7140,7220,7300,7390,7460,7520,7550,7540,7500,7470,7400,7300,7200,7090,6960,6800,6640,6500,6340,6200,6080,5950,5830,5720,5660,5610,5580,5550)
$bVals = @(45954,45954.01041666666,45954.02083333334,45954.03125,45954.04166666666,45954.05208333334,45954.0625,45954.07291666666,45954.08333333334,45954.09375,45954.10416666666,45954.11458333334,45954.125,45954.13541666666,45954.14583333334,45954.15625,45954.16666666666,45954.17708333334,45954.1875,45954.19791666666,45954.20833333334,45954.21875,45954.22916666666,45954.23958333334,45954.25,45954.26041666666,45954.27083333334,45954.28125,45954.29166666666,45954.30208333334,45954.3125,45954.32291666666,45954.33333333334,45954.34375,45954.35416666666,45954.36458333334,45954.375,45954.38541666666,45954.39583333334,45954.40625,45954.41666666666,45954.42708333334,45954.4375,45954.44791666666,45954.45833333334,45954.46875,45954.47916666666,45954.48958333334,45954.5,45954.51041666666,45954.52083333334,45954.53125,45954.54166666666,45954.55208333334,45954.5625,45954.57291666666,45954.58333333334,45954.59375,45954.60416666666,45954.61458333334,45954.625,45954.63541666666,45954.64583333334,45954.65625,45954.66666666666,45954.67708333334,45954.6875,45954.69791666666,45954.70833333334,45954.71875,45954.72916666666,45954.73958333334,45954.75,45954.76041666666,45954.77083333334,45954.78125,45954.79166666666,45954.80208333334,45954.8125,45954.82291666666,45954.83333333334,45954.84375,45954.85416666666,45954.86458333334,45954.875,45954.88541666666,45954.89583333334,45954.90625,45954.91666666666,45954.92708333334,45954.9375,45954.94791666666,45954.95833333334,45954.96875,45954.97916666666,45954.98958333334)

for ($i = 0; $i -lt $aVals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $aVals[$i]
    $ws.Cells.Item(2 + $i, 2).Value = $bVals[$i]
}
